# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for the zh-cn and de-de handback
# sheets, flips the Overview/Status text from "Ready for handoff" to
# "Handed back: in sync with en-US", and widens the columns that now
# hold longer content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cfd867b9612a30cf9b6f4809e6a6226060c4d75d/e2e/7f872013-e18f-4987-ad2e-267d4f6003b2.md"
$mdName = "7f872013-e18f-4987-ad2e-267d4f6003b2.md"

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Overview!E2 (zh-cn column) and Overview!F2 (de-de column)
if ($wsOverview.Range("E2").Value2 -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
# zh-cn!C2 and de-de!C2 (Status column)
if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) { $wsZhCn.Range("C2").Value = $newStatus }
if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) { $wsDeDe.Range("C2").Value = $newStatus }

# --- zh-cn handback info ---
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdName)
$wsZhCn.Range("J2").Value = "7f872013-e18f-4987-ad2e-267d4f6003b2.a251fc29809b721ac324c4f4d064e168aec3326a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-20 17:06:28"

# --- de-de handback info ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdName)
$wsDeDe.Range("J2").Value = "7f872013-e18f-4987-ad2e-267d4f6003b2.a251fc29809b721ac324c4f4d064e168aec3326a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-20 17:06:34"

# --- Column widths (stored xlsx width = ColumnWidth snapped to 1/6 + 5/6) ---
# Overview columns E (zh-cn) and F (de-de): widen for the longer status text
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de column C (Status): widen for the longer status text
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668

# zh-cn / de-de columns I (Latest Target File) and J (Latest Handback File):
# widen to fit the newly populated long file names
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
